# Append a new data row (row 52) to each of the four sheets, matching the
# pattern of the existing data rows (columns A-I), and let Excel extend the
# sheet's used-range / dimension automatically.

$wb = $excel.ActiveWorkbook

function Add-SensorRow($ws, $row, $timeStr, $colB, $colC, $colD, $colE, $colF, $colG, $colH, $colI) {
    $ws.Range("A$row").Value = $timeStr
    $ws.Range("B$row").Value = $colB
    $ws.Range("C$row").Value = $colC
    $ws.Range("D$row").Value = $colD
    $ws.Range("E$row").Value = $colE
    $ws.Range("F$row").Value = $colF

    # Column G holds a 24-25 digit number that must stay plain text (Excel
    # would otherwise coerce a long pure-digit string into a lossy double /
    # scientific notation). Force text formatting, assign, then strip the
    # formatting again so no stray style index is left on the cell.
    $ws.Range("G$row").NumberFormat = "@"
    $ws.Range("G$row").Value = $colG
    $ws.Range("G$row").ClearFormats()

    $ws.Range("H$row").Value = $colH
    $ws.Range("I$row").Value = $colI
}

# Sheet 1: ROW35-FE-LIFTER
$ws1 = $wb.Worksheets.Item(1)
Add-SensorRow $ws1 52 "2025-03-06 11:42:06" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c," "0x01,0x90," "0x d" 400 "568631262647113770877196" 400 13

# Sheet 2: ROW35-MID-LIFTER
$ws2 = $wb.Worksheets.Item(2)
Add-SensorRow $ws2 52 "2025-03-06 11:29:35" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x90," "0x e" 400 "568631262647113770942732" 400 14

# Sheet 3: ROW02-FE-LIFTER
$ws3 = $wb.Worksheets.Item(3)
Add-SensorRow $ws3 52 "2025-03-06 11:51:45" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c," "0x01,0x90," "0xff" 400 "568631262647113769959692" 400 255

# Sheet 4: ROW02-MID-LIFTER
$ws4 = $wb.Worksheets.Item(4)
Add-SensorRow $ws4 52 "2025-03-06 11:41:15" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x01,0x90," "0x 3" 400 "568631262647113769959692" 400 3
